$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Summary" sheet - update aggregate metrics in column B
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.96   # Current Capital
$summary.Range("B4").Value = -0.04     # Total P&L $
$summary.Range("B5").Value = -0.4      # Total P&L %
$summary.Range("B6").Value = 2         # Total Trades
$summary.Range("B8").Value = 1         # Losing Trades
$summary.Range("B9").Value = 50        # Win Rate %

# ---------------------------------------------------------------------------
# 2) "Strategy Status" sheet - update MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.95999999999999   # Capital
$status.Range("D4").Value = 2                   # Trades
$status.Range("E4").Value = -0.04                # P&L $
$status.Range("F4").Value = -0.04                # P&L %
$status.Range("G4").Value = 50                   # Win Rate %

# ---------------------------------------------------------------------------
# 3) Append the newly-closed trade (#2) as row 3 on both the "All Trades"
#    and "MarketMaking" log sheets.
# ---------------------------------------------------------------------------
function Add-TradeRow($ws) {
    $ws.Cells.Item(3, 1).Value = 2                 # Trade #

    # Force text so Excel's auto date-detection doesn't turn this into a
    # serial date number (the source data stores it as a literal string).
    $ws.Cells.Item(3, 2).NumberFormat = "@"
    $ws.Cells.Item(3, 2).Value = "2026-02-17"       # Date (kept as text)

    $ws.Cells.Item(3, 3).Value = "15:13:16"         # Time (kept as text)
    $ws.Cells.Item(3, 4).Value = "MarketMaking"     # Strategy
    $ws.Cells.Item(3, 5).Value = "UP"               # Side

    $ws.Cells.Item(3, 6).Value = 0.9                # Entry Price
    $ws.Cells.Item(3, 7).Value = 0.8100000000000001 # Exit Price

    $ws.Cells.Item(3, 8).Value = "CLOSED"           # Status

    $ws.Cells.Item(3, 9).Value = -10                # P&L %
    $ws.Cells.Item(3, 10).Value = -0.09             # P&L $
    $ws.Cells.Item(3, 11).Value = 99.95999999999999 # Capital After
    $ws.Cells.Item(3, 12).Value = 0                 # Entry Slippage (bps)
    $ws.Cells.Item(3, 13).Value = 0                 # Exit Slippage (bps)
    $ws.Cells.Item(3, 14).Value = 0.6               # Confidence

    $ws.Cells.Item(3, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item(3, 16).Value = "early_exit"      # Exit Reason

    $ws.Cells.Item(3, 17).Value = 0.11              # Duration (min)
}

Add-TradeRow $wb.Worksheets.Item("All Trades")
Add-TradeRow $wb.Worksheets.Item("MarketMaking")
